$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-22 down to 20-23.
$ws.Rows.Item(19).Insert()

# Copy the formatting (number formats, styles, etc.) from the row that was
# previously row 19 (now row 20) into the newly inserted blank row 19, so the
# new row matches the existing look (e.g. the date-formatted column D).
# Restrict the copy to the used columns (A:R) so the sheet's dimension isn't
# needlessly expanded to the full row width.
$ws.Range("A20:R20").Copy()
$ws.Range("A19:R19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44748
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100112043
$ws.Cells.Item(19, 7).Value = "Pepino dulce"
$ws.Cells.Item(19, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 250
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 18000
$ws.Cells.Item(19, 13).Value = 17500
$ws.Cells.Item(19, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 972
$ws.Cells.Item(19, 17).Value = 18
$ws.Cells.Item(19, 18).Value = "Hortaliza"

$wb.Save()
